# DAL geolocation gemaakt + getest (no coordinates)
#
# The DaGeoLocation methods (selectAll, selectOneById, selectAllByOrganism,
# insert, delete, update -> rows 31-36) had empty Auteur/Status columns.
# Fill them in the same way every other DAL class on the sheet is already
# recorded: Auteur = "Oualid", Status = "Done".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 31; $row -le 36; $row++) {
    $ws.Cells.Item($row, 3).Value = "Oualid"  # column C - Auteur
    $ws.Cells.Item($row, 4).Value = "Done"    # column D - Status
}

# Leave the same trail in the saved view state that Excel would: the
# selection moved to the row/column that was just finished editing.
$ws.Range("C33").Select()
